# Update the mislabeled decimal string in A3 (was "0.02", should read "0.03")
# and leave the selection parked on the corrected cell, matching the
# "compatible to python 2.7.x" data-fix commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "0.03"
[void]$ws.Range("A3").Select()
